$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rerun full analysis updated a few missing gene symbols (column B, "genname")
# for entrezgene IDs that previously had no matching name.
$ws.Range("B30").Value = "LRATD1"
$ws.Range("B43").Value = "SMIM3"
$ws.Range("B49").Value = "BRI3"
$ws.Range("B56").Value = "GABBR2"
